# Insert a new weekly record at row 339 (Femacal de La Calera - Ajo),
# pushing the existing rows 339:403 down to 340:404.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 339; Excel shifts rows 339:403 -> 340:404
$ws.Rows.Item(339).Insert()

# Populate the new row 339 with the new observation.
$ws.Range("A339").Value = 3
$ws.Range("B339").Value = 'Femacal de La Calera'
$ws.Range("C339").Value = 'Coquimbo'
$ws.Range("D339").Value = 44637
$ws.Range("E339").Value = 5
$ws.Range("F339").Value = 100112003
$ws.Range("G339").Value = 'Ajo'
$ws.Range("H339").Value = 'Chino'
$ws.Range("I339").Value = 'Primera'
$ws.Range("J339").Value = 85
$ws.Range("K339").Value = 16500
$ws.Range("L339").Value = 17000
$ws.Range("M339").Value = 16735
$ws.Range("N339").Value = '$/caja 10 kilos'
$ws.Range("O339").Value = 'China'
$ws.Range("P339").Value = 1674
$ws.Range("Q339").Value = 10
$ws.Range("R339").Value = 'Hortaliza'
